# Generate Report for Handback
# This script reorders / refreshes the localization-status report:
#  - Items 69621379-... and c3615a61-... have completed handback, so they move
#    to the top of each sheet with status "Handed back: in sync with en-US"
#    and updated handoff/handback timestamps.
#  - ead82ac8-... (still "In Translation") and dc988be4-... (still
#    "Ready for handoff") keep their values but shift down in the Overview
#    sheet ordering.
#  - The per-locale sheets (zh-cn, de-de) gain "Latest Target File" (F) and
#    "Latest Handback File" (G) entries (with hyperlinks) for the two items
#    that were handed back.

$wb = $excel.ActiveWorkbook

function Reset-Hyperlinks {
    param($ws)
    $existing = @()
    foreach ($hl in $ws.Hyperlinks) {
        $existing += $hl
    }
    for ($i = $existing.Count - 1; $i -ge 0; $i--) {
        $existing[$i].Delete()
    }
}

function Add-Link {
    param($ws, $cellAddr, $url, $displayText)
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, "", "", $displayText)
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "69621379-06c6-4801-b22b-fc912f443981.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value = "2016-03-23 20:19:58"

$wsOverview.Range("A3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.md"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D3").Value = "2016-03-23 20:19:58"

$wsOverview.Range("A4").Value = "ead82ac8-71de-4c3c-8905-e430a69e5b3c.md"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"
$wsOverview.Range("D4").Value = "2016-03-23 20:18:26"

$wsOverview.Range("A5").Value = "dc988be4-5b83-4861-b9e6-06e29e750867.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-23 20:19:58"

Reset-Hyperlinks $wsOverview
Add-Link $wsOverview "A2" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/69621379-06c6-4801-b22b-fc912f443981.md" "69621379-06c6-4801-b22b-fc912f443981.md"
Add-Link $wsOverview "A3" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/c3615a61-4ba7-4351-9f48-041a6b773395.md" "c3615a61-4ba7-4351-9f48-041a6b773395.md"
Add-Link $wsOverview "A4" "https://github.com/OpenLocalizationTest/oltest/blob/b9157820144708220278a7155ff95c287a323e28/e2e/ead82ac8-71de-4c3c-8905-e430a69e5b3c.md" "ead82ac8-71de-4c3c-8905-e430a69e5b3c.md"
Add-Link $wsOverview "A5" "https://github.com/OpenLocalizationTest/oltest/blob/d26aac8f2ef547d106b88cc9af7f09396f587022/e2e/dc988be4-5b83-4861-b9e6-06e29e750867.md" "dc988be4-5b83-4861-b9e6-06e29e750867.md"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "69621379-06c6-4801-b22b-fc912f443981.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-23 20:19:53"
$wsZh.Range("F2").Value = "69621379-06c6-4801-b22b-fc912f443981.md"
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Range("G2").Value = "69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.zh-cn.xlf"
$wsZh.Range("G2").Style = "HyperLink"
$wsZh.Range("H2").Value = "2016-03-23 20:20:28"
$wsZh.Range("J2").Value = "Include"

$wsZh.Range("A3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-23 20:19:53"
$wsZh.Range("F3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.md"
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Range("G3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.zh-cn.xlf"
$wsZh.Range("G3").Style = "HyperLink"
$wsZh.Range("H3").Value = "2016-03-23 20:20:28"
$wsZh.Range("J3").Value = "Include"

$wsZh.Range("A4").Value = "ead82ac8-71de-4c3c-8905-e430a69e5b3c.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "In Translation"
$wsZh.Range("D4").Value = "ead82ac8-71de-4c3c-8905-e430a69e5b3c.c04ad5d28185c320611c8629cb08172d2cf859c8.zh-cn.xlf"
$wsZh.Range("E4").Value = "2016-03-23 20:18:22"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("J4").Value = "Include"

$wsZh.Range("A5").Value = "dc988be4-5b83-4861-b9e6-06e29e750867.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "dc988be4-5b83-4861-b9e6-06e29e750867.1ac00fa429ac35b4b1acd56d8f4c929e60c73f86.zh-cn.xlf"
$wsZh.Range("E5").Value = "2016-03-23 20:19:53"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("J5").Value = "Include"

Reset-Hyperlinks $wsZh
Add-Link $wsZh "A2" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/69621379-06c6-4801-b22b-fc912f443981.md" "69621379-06c6-4801-b22b-fc912f443981.md"
Add-Link $wsZh "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e82657b987f9a12886c6489b635b494c9b736b6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.zh-cn.xlf" "69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.zh-cn.xlf"
Add-Link $wsZh "F2" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/69621379-06c6-4801-b22b-fc912f443981.md" "69621379-06c6-4801-b22b-fc912f443981.md"
Add-Link $wsZh "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e82657b987f9a12886c6489b635b494c9b736b6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.zh-cn.xlf" "69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.zh-cn.xlf"

Add-Link $wsZh "A3" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/c3615a61-4ba7-4351-9f48-041a6b773395.md" "c3615a61-4ba7-4351-9f48-041a6b773395.md"
Add-Link $wsZh "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e82657b987f9a12886c6489b635b494c9b736b6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.zh-cn.xlf" "c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.zh-cn.xlf"
Add-Link $wsZh "F3" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/c3615a61-4ba7-4351-9f48-041a6b773395.md" "c3615a61-4ba7-4351-9f48-041a6b773395.md"
Add-Link $wsZh "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e82657b987f9a12886c6489b635b494c9b736b6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.zh-cn.xlf" "c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.zh-cn.xlf"

Add-Link $wsZh "A4" "https://github.com/OpenLocalizationTest/oltest/blob/b9157820144708220278a7155ff95c287a323e28/e2e/ead82ac8-71de-4c3c-8905-e430a69e5b3c.md" "ead82ac8-71de-4c3c-8905-e430a69e5b3c.md"
Add-Link $wsZh "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50736ecc6edab5ed298f82b8e46ad3c7da5dbbea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ead82ac8-71de-4c3c-8905-e430a69e5b3c.c04ad5d28185c320611c8629cb08172d2cf859c8.zh-cn.xlf" "ead82ac8-71de-4c3c-8905-e430a69e5b3c.c04ad5d28185c320611c8629cb08172d2cf859c8.zh-cn.xlf"

Add-Link $wsZh "A5" "https://github.com/OpenLocalizationTest/oltest/blob/d26aac8f2ef547d106b88cc9af7f09396f587022/e2e/dc988be4-5b83-4861-b9e6-06e29e750867.md" "dc988be4-5b83-4861-b9e6-06e29e750867.md"
Add-Link $wsZh "D5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e82657b987f9a12886c6489b635b494c9b736b6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/dc988be4-5b83-4861-b9e6-06e29e750867.1ac00fa429ac35b4b1acd56d8f4c929e60c73f86.zh-cn.xlf" "dc988be4-5b83-4861-b9e6-06e29e750867.1ac00fa429ac35b4b1acd56d8f4c929e60c73f86.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "69621379-06c6-4801-b22b-fc912f443981.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-23 20:19:58"
$wsDe.Range("F2").Value = "69621379-06c6-4801-b22b-fc912f443981.md"
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Range("G2").Value = "69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.de-de.xlf"
$wsDe.Range("G2").Style = "HyperLink"
$wsDe.Range("H2").Value = "2016-03-23 20:20:36"
$wsDe.Range("J2").Value = "Include"

$wsDe.Range("A3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-23 20:19:58"
$wsDe.Range("F3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.md"
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Range("G3").Value = "c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.de-de.xlf"
$wsDe.Range("G3").Style = "HyperLink"
$wsDe.Range("H3").Value = "2016-03-23 20:20:36"
$wsDe.Range("J3").Value = "Include"

$wsDe.Range("A4").Value = "ead82ac8-71de-4c3c-8905-e430a69e5b3c.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "In Translation"
$wsDe.Range("D4").Value = "ead82ac8-71de-4c3c-8905-e430a69e5b3c.c04ad5d28185c320611c8629cb08172d2cf859c8.de-de.xlf"
$wsDe.Range("E4").Value = "2016-03-23 20:18:26"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("J4").Value = "Include"

$wsDe.Range("A5").Value = "dc988be4-5b83-4861-b9e6-06e29e750867.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "dc988be4-5b83-4861-b9e6-06e29e750867.1ac00fa429ac35b4b1acd56d8f4c929e60c73f86.de-de.xlf"
$wsDe.Range("E5").Value = "2016-03-23 20:19:58"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("J5").Value = "Include"

Reset-Hyperlinks $wsDe
Add-Link $wsDe "A2" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/69621379-06c6-4801-b22b-fc912f443981.md" "69621379-06c6-4801-b22b-fc912f443981.md"
Add-Link $wsDe "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f944c176c4413ca5c1c7095645bbd6807661ed4f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.de-de.xlf" "69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.de-de.xlf"
Add-Link $wsDe "F2" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/69621379-06c6-4801-b22b-fc912f443981.md" "69621379-06c6-4801-b22b-fc912f443981.md"
Add-Link $wsDe "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f944c176c4413ca5c1c7095645bbd6807661ed4f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.de-de.xlf" "69621379-06c6-4801-b22b-fc912f443981.5257f64c79751262054a573e719fd2dd387b3fa0.de-de.xlf"

Add-Link $wsDe "A3" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/c3615a61-4ba7-4351-9f48-041a6b773395.md" "c3615a61-4ba7-4351-9f48-041a6b773395.md"
Add-Link $wsDe "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f944c176c4413ca5c1c7095645bbd6807661ed4f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.de-de.xlf" "c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.de-de.xlf"
Add-Link $wsDe "F3" "https://github.com/OpenLocalizationTest/oltest/blob/0353f64a738b4d578dcaa8e50e0284688a33cbc4/e2e/c3615a61-4ba7-4351-9f48-041a6b773395.md" "c3615a61-4ba7-4351-9f48-041a6b773395.md"
Add-Link $wsDe "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f944c176c4413ca5c1c7095645bbd6807661ed4f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.de-de.xlf" "c3615a61-4ba7-4351-9f48-041a6b773395.86fa06e7d30693d167a72b8433021e3af7c38e60.de-de.xlf"

Add-Link $wsDe "A4" "https://github.com/OpenLocalizationTest/oltest/blob/b9157820144708220278a7155ff95c287a323e28/e2e/ead82ac8-71de-4c3c-8905-e430a69e5b3c.md" "ead82ac8-71de-4c3c-8905-e430a69e5b3c.md"
Add-Link $wsDe "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e62c9a320481d6100d0b4f14e0b29e9437fd3f7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ead82ac8-71de-4c3c-8905-e430a69e5b3c.c04ad5d28185c320611c8629cb08172d2cf859c8.de-de.xlf" "ead82ac8-71de-4c3c-8905-e430a69e5b3c.c04ad5d28185c320611c8629cb08172d2cf859c8.de-de.xlf"

Add-Link $wsDe "A5" "https://github.com/OpenLocalizationTest/oltest/blob/d26aac8f2ef547d106b88cc9af7f09396f587022/e2e/dc988be4-5b83-4861-b9e6-06e29e750867.md" "dc988be4-5b83-4861-b9e6-06e29e750867.md"
Add-Link $wsDe "D5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f944c176c4413ca5c1c7095645bbd6807661ed4f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/dc988be4-5b83-4861-b9e6-06e29e750867.1ac00fa429ac35b4b1acd56d8f4c929e60c73f86.de-de.xlf" "dc988be4-5b83-4861-b9e6-06e29e750867.1ac00fa429ac35b4b1acd56d8f4c929e60c73f86.de-de.xlf"
